# Update the Ylo column (D) values for the image rows from 1 to 200.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D6").Value = 200

# Move the active selection to D5 (as left by the author after editing).
$ws.Range("D5").Select()
